# Re-grade student 48966 after quality standards README detection
# Row 9 corresponds to Student ID 48966 on the "Grade Comparison" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grade Comparison")

# B9: Calculated Grade (10 Skills) 79 -> 80.5
$ws.Range("B9").Value = 80.5

# C9: Weighted Grade (Final) 64.91220900398571 -> 68.03772334967967
$ws.Range("C9").Value = 68.03772334967967

# D9: Difference -14.08779099601429 -> -12.46227665032033
$ws.Range("D9").Value = -12.46227665032033

# F9: Penalty 14.08779099601428 -> 12.46227665032033
$ws.Range("F9").Value = 12.46227665032033

# G9: Performance Tier "Potential" -> "Good"
$ws.Range("G9").Value = "Good"
